$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 28 }

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $val -ne "") {
            $trimmed = $val.TrimStart()
            $encoded = $trimmed.Replace(" ", "%20")
            $cell.Value = $encoded
        }
    }
}
